# maj suite à formation
# Update the contact e-mail address shown on the title slide's subtitle
# placeholder: "Philippe.Renevier@ac-grenoble.fr" -> "Philippe.Renevier-Gonin@ac-grenoble.fr"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$oldEmail = "Philippe.Renevier@ac-grenoble.fr"
$newEmail = "Philippe.Renevier-Gonin@ac-grenoble.fr"

$idx = $tr.Text.IndexOf($oldEmail)
if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, $oldEmail.Length)
    $target.Text = $newEmail
}
